$wb = $excel.ActiveWorkbook

# --- Sheet "component" (sheet2) ---
$ws2 = $wb.Worksheets.Item("component")

# "close" -> "equal" relabel, and "No`nuse" -> "no`nuse" casing fix.
$ws2.Range("I1").Value = "equal"
$ws2.Range("F1").Value = "no`nuse"

[void]$ws2.Range("J2").Select()

# --- Sheet "parameter" (sheet1) ---
$ws1 = $wb.Worksheets.Item("parameter")
$ws1.Activate()

# The sheet currently has an unused blank row 1 (data starts at row 2).
# Remove it so the header moves from row 2 to row 1 and the data rows
# shift up from 3..7 to 2..6.
$ws1.Rows.Item(1).Delete()

# New header cell + relabel the "close" threshold column as "equal".
$ws1.Range("B1").Value = "Parameter name"
$ws1.Range("F1").Value = "equal"

# Update selection to match the edited workbook (this sheet stays active).
[void]$ws1.Range("B2").Select()
